$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.487.04'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.60%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.353.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.05%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '190.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '559.45'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.344.39'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.93%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.584'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.32%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.182'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.587'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.26'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000271'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.878.17'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.98%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '605.07'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.508.92'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.64%  '

$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.86%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.118'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.331.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.70%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.906'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.58%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '100.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.97%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.75'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.82'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.91'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '583.53'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.45%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.09'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.53%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.105'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.70%  '

$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.713.14'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.43%  '

$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '57.20'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.33%  '

$ws.Range("B39").Value = 'Dai'
$ws.Range("C39").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.21%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.55'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '33.92'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.78%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.91%  '

$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.129'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.51%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0710'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.79%  '

$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.69'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.81%  '

$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.39'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.88%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.341'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0421'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.130'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.46%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.58'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.34%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.08%  '
